# Update NATMI edge-weight stats for Pros1-Mertk following Dr Hou's advice
# (ligand/receptor-expressing-cell counts go from 1 to 3, and all dependent
#  totals / specificities are recomputed accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: row number, followed by (column index, new value) pairs
# for every cell that the new pipeline run touched.
$updates = @(
    @(2,5,3,7,28.48226033333333,8,85.44678099999999,9,0.2101651977164657,10,0.2101651977164658,11,3,13,11.76011933333333,14,35.280358,15,0.1622550713151862,16,0.1622550713151862,17,334.9547804030664,18,3014.593023627598,19,0.03410036914345536,20,0.03410036914345536),
    @(3,5,3,7,28.48226033333333,8,85.44678099999999,9,0.2101651977164657,10,0.2101651977164658,11,3,13,3.903127666666666,14,11.709383,15,0.05385168636105758,16,0.05385168636105758,17,111.1698983162359,18,1000.529084846123,19,0.01131775031143677,20,0.01131775031143677),
    @(4,5,3,7,28.48226033333333,8,85.44678099999999,9,0.2101651977164657,10,0.2101651977164658,11,3,13,50.38343933333334,14,151.150318,15,0.6951433323438234,16,0.6951433323438234,17,1435.034235580706,18,12915.30812022636,19,0.1460949358833225,20,0.1460949358833225),
    @(5,5,3,7,28.48226033333333,8,85.44678099999999,9,0.2101651977164657,10,0.2101651977164658,11,3,13,6.432523333333333,14,19.29757,15,0.08874990997993266,16,0.08874990997993268,17,183.2128041802411,18,1648.91523762217,19,0.01865214237825108,20,0.01865214237825109),
    @(6,5,3,7,83.45109033333334,8,250.353271,9,0.6157697701763504,10,0.6157697701763504,11,3,13,11.76011933333333,14,35.280358,15,0.1622550713151862,16,0.1622550713151862,17,981.3947808167799,18,8832.553027351018,19,0.09991176797369954,20,0.09991176797369954),
    @(7,5,3,7,83.45109033333334,8,250.353271,9,0.6157697701763504,10,0.6157697701763504,11,3,13,3.903127666666666,14,11.709383,15,0.05385168636105758,16,0.05385168636105758,17,325.7202594935326,18,2931.482335441793,19,0.03316024053415732,20,0.03316024053415733),
    @(8,5,3,7,83.45109033333334,8,250.353271,9,0.6157697701763504,10,0.6157697701763504,11,3,13,50.38343933333334,14,151.150318,15,0.6951433323438234,16,0.6951433323438234,17,4204.55294711002,18,37840.97652399018,19,0.4280482499969785,20,0.4280482499969785),
    @(9,5,3,7,83.45109033333334,8,250.353271,9,0.6157697701763504,10,0.6157697701763504,11,3,13,6.432523333333333,14,19.29757,15,0.08874990997993266,16,0.08874990997993268,17,536.8010857612745,18,4831.209771851471,19,0.05464951167151492,20,0.05464951167151493),
    @(10,5,3,7,14.78130366666667,8,44.343911,9,0.1090684366779874,10,0.1090684366779875,11,3,13,11.76011933333333,14,35.280358,15,0.1622550713151862,16,0.1622550713151862,17,173.8298950222376,18,1564.469055200138,19,0.01769690697142272,20,0.01769690697142273),
    @(11,5,3,7,14.78130366666667,8,44.343911,9,0.1090684366779874,10,0.1090684366779875,11,3,13,3.903127666666666,14,11.709383,15,0.05385168636105758,16,0.05385168636105758,17,57.6933152907681,18,519.239837616913,19,0.005873519243873849,20,0.00587351924387385),
    @(12,5,3,7,14.78130366666667,8,44.343911,9,0.1090684366779874,10,0.1090684366779875,11,3,13,50.38343933333334,14,151.150318,15,0.6951433323438234,16,0.6951433323438234,17,744.7329165570776,18,6702.596249013698,19,0.07581819652586749,20,0.07581819652586749),
    @(13,5,3,7,14.78130366666667,8,44.343911,9,0.1090684366779874,10,0.1090684366779875,11,3,13,6.432523333333333,14,19.29757,15,0.08874990997993266,16,0.08874990997993268,17,95.08108073291888,18,855.72972659627,19,0.009679813936823373,20,0.009679813936823375),
    @(14,5,3,7,8.808546666666667,8,26.42564,9,0.0649965954291964,10,0.06499659542919642,11,3,13,11.76011933333333,14,35.280358,15,0.1622550713151862,16,0.1622550713151862,17,103.5895599532356,18,932.30603957912,19,0.01054602722660857,20,0.01054602722660857),
    @(15,5,3,7,8.808546666666667,8,26.42564,9,0.0649965954291964,10,0.06499659542919642,11,3,13,3.903127666666666,14,11.709383,15,0.05385168636105758,16,0.05385168636105758,17,34.38088219779111,18,309.42793978012,19,0.003500176271589633,20,0.003500176271589634),
    @(16,5,3,7,8.808546666666667,8,26.42564,9,0.0649965954291964,10,0.06499659542919642,11,3,13,50.38343933333334,14,151.150318,15,0.6951433323438234,16,0.6951433323438234,17,443.8048765948355,18,3994.24388935352,19,0.04518194993765491,20,0.04518194993765492),
    @(17,5,3,7,8.808546666666667,8,26.42564,9,0.0649965954291964,10,0.06499659542919642,11,3,13,6.432523333333333,14,19.29757,15,0.08874990997993266,16,0.08874990997993268,17,56.66118196608889,18,509.9506376948,19,0.005768441993343284,20,0.005768441993343286)
)

foreach ($rowUpdate in $updates) {
    $r = $rowUpdate[0]
    for ($i = 1; $i -lt $rowUpdate.Length; $i += 2) {
        $col = $rowUpdate[$i]
        $val = $rowUpdate[$i + 1]
        $ws.Cells.Item($r, $col).Value = $val
    }
}